# Restored from revision of admin on 01/06/2021 08:29:02 AM.TEST Author: admin. Type: SAVE.
# Change: cell C10 on the active sheet changes its numeric value from 18 to 1
# (the cell's style/formatting, s="20", is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
